# Update the EPEX Spot public price workbook with the 06-jul data point.
#
# 1) "Prix Spot" sheet: append a new day column (W) "06-jul" with its 24
#    hourly prices.
# 2) "Gaz" sheet: append a new row (20) for 2025-07-04.
# 3) "CO2" sheet: append a new row (20) for 2025-07-04.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add column W = "06-jul"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold font, border, centered
# alignment) from the previous day's header (V1) onto the new one (W1),
# then set its text.
$ws1.Range("V1").Copy($ws1.Range("W1"))
$ws1.Range("W1").Value = "06-jul"

$hourlyPrices = @{
    2  = 86.18000000000001
    3  = 72.55
    4  = 57.04
    5  = 41.43
    6  = 39.28
    7  = 36.98
    8  = 33.94
    9  = 32.11
    10 = 35.35
    11 = 11.14
    12 = 18.49
    13 = 14.94
    14 = 12.69
    15 = 9.58
    16 = 2.99
    17 = 0.65
    18 = 0.65
    19 = 0.65
    20 = 11.22
    21 = 31.53
    22 = 31.61
    23 = 35.61
    24 = 70.19
    25 = 71.44
}

foreach ($row in $hourlyPrices.Keys) {
    # Column W is the 23rd column.
    $ws1.Cells.Item($row, 23).Value = $hourlyPrices[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": add row 20 = 2025-07-04 / 32.775
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date to be stored as plain text (matching the existing
# A2:A19 cells) instead of letting Excel auto-convert it to a date
# serial number, then drop back to the default "Normal" style so the
# new row matches the unstyled look of the existing data rows.
$ws2.Range("A20").NumberFormat = "@"
$ws2.Range("A20").Value = "2025-07-04"
$ws2.Range("A20").Style = "Normal"
$ws2.Range("B20").Value = 32.775

# ---------------------------------------------------------------------
# Sheet "CO2": add row 20 = 2025-07-04 / 70.92
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A20").NumberFormat = "@"
$ws3.Range("A20").Value = "2025-07-04"
$ws3.Range("A20").Style = "Normal"
$ws3.Range("B20").Value = 70.92
